$d = $word.ActiveDocument

$pairs = @(
    @("766÷5=", "420÷2="),
    @("958÷5=", "687÷8="),
    @("600÷9=", "240÷6="),
    @("335÷6=", "623÷6="),
    @("595÷3=", "391÷6="),
    @("978÷8=", "253÷2="),
    @("582÷2=", "719÷9="),
    @("129÷4=", "317÷3="),
    @("526÷7=", "530÷8="),
    @("780÷7=", "859÷9="),
    @("441÷3=", "803÷2="),
    @("717÷2=", "969÷8="),
    @("980÷6=", "157÷3="),
    @("302÷7=", "658÷6="),
    @("282÷2=", "880÷8="),
    @("984÷6=", "271÷8="),
    @("820÷5=", "265÷8="),
    @("759÷8=", "682÷9="),
    @("157÷9=", "155÷3="),
    @("875÷3=", "672÷5="),
    @("188÷2=", "317÷6="),
    @("563÷9=", "544÷8="),
    @("489÷6=", "352÷8="),
    @("710÷5=", "215÷2="),
    @("634÷9=", "408÷8=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
